$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C6").Value = 2.890185992428072
$ws.Range("D6").Value = 3.714905814451302
$ws.Range("E6").Value = 0.4179997037158225
$ws.Range("F6").Value = 1.543688107481775
$ws.Range("G6").Value = 1.72069282810891
$ws.Range("C7").Value = 4.016309467720443
$ws.Range("D7").Value = 0.1300111556989939
$ws.Range("E7").Value = 0.4901371972304145
$ws.Range("F7").Value = 2.049247970961359
$ws.Range("G7").Value = 1.072177638005135
$ws.Range("C8").Value = 1.450728809180919
$ws.Range("D8").Value = 2.918115882885407
$ws.Range("E8").Value = 0.6100846378576197
$ws.Range("F8").Value = 1.229658311504014
$ws.Range("G8").Value = 0.8051804048969398
$ws.Range("C9").Value = 2.362842307729407
$ws.Range("D9").Value = 5.052065509201971
$ws.Range("E9").Value = 1.308846612043794
$ws.Range("F9").Value = 1.514483439887546
$ws.Range("G9").Value = 2.156965655612491
$ws.Range("C10").Value = 3.600943004536212
$ws.Range("D10").Value = 4.400904943103417
$ws.Range("E10").Value = 1.404581544613537
$ws.Range("F10").Value = 1.465089618404473
$ws.Range("G10").Value = 1.859439426993244
$ws.Range("C11").Value = 4.895272577004197
$ws.Range("D11").Value = 8.380442961710859
$ws.Range("E11").Value = 4.115704686220278
$ws.Range("F11").Value = 1.534052611338125
$ws.Range("G11").Value = 1.807094929033588
$ws.Range("C12").Value = 3.556280773819503
$ws.Range("D12").Value = 6.743409376906882
$ws.Range("E12").Value = 5.339025292833842
$ws.Range("F12").Value = 1.761070506569439
$ws.Range("G12").Value = 2.472289524208387
$ws.Range("C13").Value = 6.53482650340761
$ws.Range("D13").Value = 5.77785891826038
$ws.Range("E13").Value = 4.168688394382894
$ws.Range("F13").Value = 1.716797499598147
$ws.Range("G13").Value = 2.702524820403045
$ws.Range("C14").Value = 5.216467899905491
$ws.Range("D14").Value = 8.147669858029859
$ws.Range("E14").Value = 5.402905851733751
$ws.Range("F14").Value = 1.83533289546119
$ws.Range("G14").Value = 2.604037549316263
$ws.Range("C15").Value = 5.161883454294133
$ws.Range("D15").Value = 8.721037702202484
$ws.Range("E15").Value = 5.540399034816313
$ws.Range("F15").Value = 1.990818339451958
$ws.Range("G15").Value = 2.847442186286497
$ws.Range("C16").Value = 5.772280344618859
$ws.Range("D16").Value = 9.178491449275816
$ws.Range("E16").Value = 5.483974957733746
$ws.Range("F16").Value = 1.802887611378982
$ws.Range("G16").Value = 4.175350867531982
$ws.Range("C17").Value = 4.653627239215556
$ws.Range("D17").Value = 1.030705509951482
$ws.Range("E17").Value = 4.066951309461864
$ws.Range("F17").Value = 2.195359635464623
$ws.Range("G17").Value = 4.046325449774004
$ws.Range("C18").Value = 4.531042037695504
$ws.Range("D18").Value = 5.480827786999743
$ws.Range("E18").Value = 3.498728062791118
$ws.Range("F18").Value = 2.282756903655365
$ws.Range("G18").Value = 2.774640648251451
$ws.Range("C19").Value = 4.761101999744116
$ws.Range("D19").Value = 4.363343947245058
$ws.Range("E19").Value = 2.165899987081528
$ws.Range("F19").Value = 2.396366066936376
$ws.Range("G19").Value = 2.298267337383058
